# Update the CNPJ tracking sheet:
#  - Row 2: replace the Gustavo Xavier Maurmann record with a new
#    "AGENCIA DE DADOS LTDA" record (new CNPJ + full address/contact info).
#  - Row 3: unchanged (iFood record).
#  - Row 4: the Kelvin Crisostomo Gomes record is replaced by "N/A"
#    placeholders across every column.
#  - Row 5: becomes the (previously-row-2) Gustavo Xavier Maurmann record;
#    it no longer has a FANTASIA (C) value.
#  - Row 6 (new): Lucas Caleb Alves Peixoto de Quadro record, mirroring the
#    Gustavo row's address/municipality but with its own CNPJ/contact info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: AGENCIA DE DADOS LTDA ---------------------------------------
$ws.Range("A2").Value = 29148959000150
$ws.Range("B2").Value = "AGENCIA DE DADOS LTDA"
$ws.Range("C2").Value = "AGENCIA DE DADOS"
$ws.Range("D2").Value = "MICRO EMPRESA"
$ws.Range("E2").Value = "SETOR COMERCIAL SUL QUADRA 01 BLOCO G"
$ws.Range("F2").Value = "BRASILIA"
$ws.Range("G2").Value = "ASA SUL"
$ws.Range("H2").Value = "DF"
$ws.Range("I2").Value = "70.309-900"
$ws.Range("J2").Value = "administrativo@intraseg.com.br"
$ws.Range("K2").Value = "(61) 8549-3011"

# --- Row 3: unchanged ----------------------------------------------------

# --- Row 4: Kelvin record becomes an all "N/A" row ----------------------
$ws.Range("B4").Value = "N/A"
$ws.Range("C4").Value = "N/A"
$ws.Range("D4").Value = "N/A"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = "N/A"
$ws.Range("H4").Value = "N/A"
$ws.Range("I4").Value = "N/A"
$ws.Range("J4").Value = "N/A"
$ws.Range("K4").Value = "N/A"

# --- Row 5: becomes the Gustavo Xavier Maurmann record -------------------
$ws.Range("C5").ClearContents()
$ws.Range("B5").Value = "53.405.671 GUSTAVO XAVIER MAURMANN"
$ws.Range("D5").Value = "MICRO EMPRESA"
$ws.Range("E5").Value = "QUADRA CNA 3"
$ws.Range("F5").Value = "BRASILIA"
$ws.Range("G5").Value = "TAGUATINGA NORTE (TAGUATINGA)"
$ws.Range("H5").Value = "DF"
$ws.Range("I5").Value = "72.110-035"
$ws.Range("J5").Value = "gustavo.maurmann2@gmail.com"
$ws.Range("K5").Value = "(61) 9525-3612"

# --- Row 6 (new): Lucas Caleb Alves Peixoto de Quadro --------------------
$ws.Range("A6").Value = 51041667000173
$ws.Range("B6").Value = "51.041.667 LUCAS CALEB ALVES PEIXOTO DE QUADRO"
$ws.Range("D6").Value = "MICRO EMPRESA"
$ws.Range("E6").Value = "QUADRA CNA 3"
$ws.Range("F6").Value = "BRASILIA"
$ws.Range("G6").Value = "TAGUATINGA NORTE (TAGUATINGA)"
$ws.Range("H6").Value = "DF"
$ws.Range("I6").Value = "72.110-035"
$ws.Range("J6").Value = "soucaleb51@gmail.com"
$ws.Range("K6").Value = "(66) 9156-4931"
